$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.950.79'
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').Value = '1.758.81'
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('D4').Value = '''1.004'
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').Value = '''335.85'
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('D7').Value = '''0.3828'
$ws.Range('E7').Value = '  -2.17%  '
$ws.Range('D8').Value = '''0.3386'
$ws.Range('E8').Value = '  -2.18%  '
$ws.Range('D9').Value = '''44.85'
$ws.Range('E9').Value = '  -7.16%  '
$ws.Range('D10').Value = '''1.113'
$ws.Range('E10').Value = '  -4.79%  '
$ws.Range('D11').Value = '''0.07221'
$ws.Range('E11').Value = '  -4.55%  '
$ws.Range('E12').Value = '  -0.30%  '
$ws.Range('D13').Value = '''22.29'
$ws.Range('E13').Value = '  -1.57%  '
$ws.Range('D14').Value = '''6.134'
$ws.Range('E14').Value = '  -5.49%  '
$ws.Range('D15').Value = '''7.170'
$ws.Range('E15').Value = '  -0.63%  '
$ws.Range('D16').Value = '1.758.00'
$ws.Range('E16').Value = '  -1.82%  '
$ws.Range('D17').Value = '''0.00001057'
$ws.Range('E17').Value = '  -3.70%  '
$ws.Range('D18').Value = '''0.06610'
$ws.Range('E18').Value = '  -2.04%  '
$ws.Range('D19').Value = '''79.25'
$ws.Range('E19').Value = '  -5.92%  '
$ws.Range('D20').Value = '''1.001'
$ws.Range('E20').Value = '  -0.32%  '
$ws.Range('D21').Value = '''16.63'
$ws.Range('E21').Value = '  -6.62%  '
$ws.Range('D22').Value = '''6.220'
$ws.Range('E22').Value = '  -5.55%  '
$ws.Range('D23').Value = '27.971.47'
$ws.Range('E23').Value = '  +0.54%  '
$ws.Range('D24').Value = '''11.62'
$ws.Range('E24').Value = '  -6.78%  '
$ws.Range('D25').Value = '''2.391'
$ws.Range('E25').Value = '  -0.73%  '
$ws.Range('D26').Value = '''152.31'
$ws.Range('E26').Value = '  -2.27%  '
$ws.Range('D27').Value = '''19.77'
$ws.Range('E27').Value = '  -7.18%  '
$ws.Range('D28').Value = '''2.307'
$ws.Range('E28').Value = '  -8.64%  '
$ws.Range('D29').Value = '1.958.09'
$ws.Range('E29').Value = '  -1.71%  '
$ws.Range('D30').Value = '''1.267'
$ws.Range('E30').Value = '  -16.34%  '
$ws.Range('D31').Value = '''131.63'
$ws.Range('E31').Value = '  -4.53%  '
$ws.Range('D32').Value = '''4.013'
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('D33').Value = '''5.831'
$ws.Range('E33').Value = '  -8.27%  '
$ws.Range('D34').Value = '''0.08817'
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('E35').Value = '  -7.60%  '
$ws.Range('D36').Value = '''0.6599'
$ws.Range('E36').Value = '  -5.31%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '''0.06220'
$ws.Range('E37').Value = '  -4.63%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '''5.166'
$ws.Range('E38').Value = '  -6.41%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.02289'
$ws.Range('E39').Value = '  -7.61%  '
$ws.Range('D40').Value = '''0.2110'
$ws.Range('E40').Value = '  -6.15%  '
$ws.Range('B41').Value = 'WEMIXTOKEN'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').Value = '''1.477'
$ws.Range('E41').Value = '  -6.50%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '''1.208'
$ws.Range('E42').Value = '  -4.92%  '
$ws.Range('D43').Value = '''7.961'
$ws.Range('E43').Value = '  -7.60%  '
$ws.Range('D44').Value = '''1.000'
$ws.Range('E44').Value = '  -0.34%  '
$ws.Range('D45').Value = '''13.83'
$ws.Range('E45').Value = '  -6.34%  '
$ws.Range('D46').Value = '''3.816'
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('D47').Value = '''0.6028'
$ws.Range('E47').Value = '  -6.18%  '
$ws.Range('D48').Value = '''126.25'
$ws.Range('E48').Value = '  -5.88%  '
$ws.Range('D49').Value = '''2.007'
$ws.Range('E49').Value = '  -7.40%  '
$ws.Range('B50').Value = 'Flow'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range('D50').Value = '''1.117'
$ws.Range('E50').Value = '  +2.45%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').Value = '''1.172'
$ws.Range('E51').Value = '  -0.38%  '
